$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (LSR no.) to Text format for new rows so numeric-looking values stay text
$ws.Range("C675:C695").NumberFormat = "@"

# Row 675
$ws.Range("C675").Value = "1"
$ws.Range("E675").Value = "nausea_n  "
$ws.Range("F675").Value = "Number of participants with measurement for nausea following intervention"
$ws.Range("G675").Value = "GMHO:0000183"
$ws.Range("H675").Value = "nausea following intervention"
$ws.Range("I675").Value = "A digestive system adverse event following intervention that has an outcome of nausea, a gastric discomfort associated with the inclination to vomit"
$ws.Range("J675").Value = "digestive system adverse event following intervention"
$ws.Range("K675").Value = "Intervention outcomes and spillover effects"
$ws.Range("L675").Value = "GMHO:0000178"
$ws.Range("M675").Value = "COMBO"
$ws.Range("N675").Value = "GMHO:0000183,GMHO:0000204"

# Row 676
$ws.Range("C676").Value = "1"
$ws.Range("E676").Value = "nausea_n  "
$ws.Range("F676").Value = "Number of participants with measurement for nausea following intervention"
$ws.Range("G676").Value = "GMHO:0000206"
$ws.Range("H676").Value = "number of participants with measurement"
$ws.Range("I676").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J676").Value = "number of intervention participants"
$ws.Range("K676").Value = "Intervention outcomes and spillover effects"
$ws.Range("L676").Value = "GMHO:0000183,GMHO:0000204"
$ws.Range("M676").Value = "No Combo"

# Row 677
$ws.Range("C677").Value = "1"
$ws.Range("E677").Value = "nausea_n  "
$ws.Range("F677").Value = "Number of participants with measurement for nausea following intervention"
$ws.Range("G677").Value = "GMHO:0000204"
$ws.Range("H677").Value = "measurement datum at followup"
$ws.Range("I677").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J677").Value = "measurement datum"
$ws.Range("K677").Value = "Intervention outcomes and spillover effects"
$ws.Range("L677").Value = "GMHO:0000178"
$ws.Range("M677").Value = "COMBO"
$ws.Range("N677").Value = "GMHO:0000183,GMHO:0000204"

# Row 678
$ws.Range("C678").Value = "1"
$ws.Range("E678").Value = "headache_n"
$ws.Range("F678").Value = "Number of participants with measurement for headache following intervention"
$ws.Range("G678").Value = "GMHO:0000181"
$ws.Range("H678").Value = "headache following intervention"
$ws.Range("I678").Value = "A pain following intervention that has an outcome of headache."
$ws.Range("J678").Value = "pain following intervention"
$ws.Range("K678").Value = "Intervention outcomes and spillover effects"
$ws.Range("L678").Value = "GMHO:0000184"
$ws.Range("M678").Value = "COMBO"
$ws.Range("N678").Value = "GMHO:0000181,GMHO:0000204"

# Row 679
$ws.Range("C679").Value = "1"
$ws.Range("E679").Value = "headache_n"
$ws.Range("F679").Value = "Number of participants with measurement for headache following intervention"
$ws.Range("G679").Value = "GMHO:0000206"
$ws.Range("H679").Value = "number of participants with measurement"
$ws.Range("I679").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J679").Value = "number of intervention participants"
$ws.Range("K679").Value = "Intervention outcomes and spillover effects"
$ws.Range("L679").Value = "GMHO:0000181,GMHO:0000204"
$ws.Range("M679").Value = "No Combo"

# Row 680
$ws.Range("C680").Value = "1"
$ws.Range("E680").Value = "headache_n"
$ws.Range("F680").Value = "Number of participants with measurement for headache following intervention"
$ws.Range("G680").Value = "GMHO:0000204"
$ws.Range("H680").Value = "measurement datum at followup"
$ws.Range("I680").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J680").Value = "measurement datum"
$ws.Range("K680").Value = "Intervention outcomes and spillover effects"
$ws.Range("L680").Value = "GMHO:0000184"
$ws.Range("M680").Value = "COMBO"
$ws.Range("N680").Value = "GMHO:0000181,GMHO:0000204"

# Row 681
$ws.Range("C681").Value = "1"
$ws.Range("E681").Value = "insomnia_n"
$ws.Range("F681").Value = "Number of participants with measurement for insomnia following intervention"
$ws.Range("G681").Value = "GMHO:0000182"
$ws.Range("H681").Value = "insomnia following intervention"
$ws.Range("I681").Value = "Recuperando datos. Espere unos segundos e intente cortar o copiar de nuevo."
$ws.Range("J681").Value = "adverse event following an intervention"
$ws.Range("K681").Value = "Intervention outcomes and spillover effects"
$ws.Range("L681").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M681").Value = "COMBO"
$ws.Range("N681").Value = "GMHO:0000182,GMHO:0000204"

# Row 682
$ws.Range("C682").Value = "1"
$ws.Range("E682").Value = "insomnia_n"
$ws.Range("F682").Value = "Number of participants with measurement for insomnia following intervention"
$ws.Range("G682").Value = "GMHO:0000206"
$ws.Range("H682").Value = "number of participants with measurement"
$ws.Range("I682").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J682").Value = "number of intervention participants"
$ws.Range("K682").Value = "Intervention outcomes and spillover effects"
$ws.Range("L682").Value = "GMHO:0000182,GMHO:0000204"
$ws.Range("M682").Value = "No Combo"

# Row 683
$ws.Range("C683").Value = "1"
$ws.Range("E683").Value = "insomnia_n"
$ws.Range("F683").Value = "Number of participants with measurement for insomnia following intervention"
$ws.Range("G683").Value = "GMHO:0000204"
$ws.Range("H683").Value = "measurement datum at followup"
$ws.Range("I683").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J683").Value = "measurement datum"
$ws.Range("K683").Value = "Intervention outcomes and spillover effects"
$ws.Range("L683").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M683").Value = "COMBO"
$ws.Range("N683").Value = "GMHO:0000182,GMHO:0000204"

# Row 684
$ws.Range("C684").Value = "1"
$ws.Range("E684").Value = "dry_mouth_n"
$ws.Range("F684").Value = "Number of participants with measurement for dry mouth following intervention"
$ws.Range("G684").Value = "GMHO:0000180"
$ws.Range("H684").Value = "dry mouth following intervention"
$ws.Range("I684").Value = "A digestive system adverse event following intervention that involves experiencing dry mouth."
$ws.Range("J684").Value = "digestive system adverse event following intervention"
$ws.Range("K684").Value = "Intervention outcomes and spillover effects"
$ws.Range("L684").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M684").Value = "COMBO"
$ws.Range("N684").Value = "GMHO:0000180,GMHO:0000204"

# Row 685
$ws.Range("C685").Value = "1"
$ws.Range("E685").Value = "dry_mouth_n"
$ws.Range("F685").Value = "Number of participants with measurement for dry mouth following intervention"
$ws.Range("G685").Value = "GMHO:0000206"
$ws.Range("H685").Value = "number of participants with measurement"
$ws.Range("I685").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J685").Value = "number of intervention participants"
$ws.Range("K685").Value = "Intervention outcomes and spillover effects"
$ws.Range("L685").Value = "GMHO:0000180,GMHO:0000204"
$ws.Range("M685").Value = "No Combo"

# Row 686
$ws.Range("C686").Value = "1"
$ws.Range("E686").Value = "dry_mouth_n"
$ws.Range("F686").Value = "Number of participants with measurement for dry mouth following intervention"
$ws.Range("G686").Value = "GMHO:0000204"
$ws.Range("H686").Value = "measurement datum at followup"
$ws.Range("I686").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J686").Value = "measurement datum"
$ws.Range("K686").Value = "Intervention outcomes and spillover effects"
$ws.Range("L686").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M686").Value = "COMBO"
$ws.Range("N686").Value = "GMHO:0000180,GMHO:0000204"

# Row 687
$ws.Range("C687").Value = "1"
$ws.Range("E687").Value = "constipation_n"
$ws.Range("F687").Value = "Number of participants with measurement for constipation following intervention"
$ws.Range("G687").Value = "GMHO:0000257"
$ws.Range("H687").Value = "constipation adverse event following intervention"
$ws.Range("I687").Value = "A digestive system adverse event following an intervention that involves the abnormally delayed or infrequent passage of dry hardened feces."
$ws.Range("J687").Value = "digestive system adverse event following intervention"
$ws.Range("K687").Value = "Intervention outcomes and spillover effects"
$ws.Range("L687").Value = "GMHO:0000178"
$ws.Range("M687").Value = "COMBO"
$ws.Range("N687").Value = "GMHO:0000257,GMHO:0000204"

# Row 688
$ws.Range("C688").Value = "1"
$ws.Range("E688").Value = "constipation_n"
$ws.Range("F688").Value = "Number of participants with measurement for constipation following intervention"
$ws.Range("G688").Value = "GMHO:0000206"
$ws.Range("H688").Value = "number of participants with measurement"
$ws.Range("I688").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J688").Value = "number of intervention participants"
$ws.Range("K688").Value = "Intervention outcomes and spillover effects"
$ws.Range("L688").Value = "GMHO:0000257,GMHO:0000204"
$ws.Range("M688").Value = "No Combo"

# Row 689
$ws.Range("C689").Value = "1"
$ws.Range("E689").Value = "constipation_n"
$ws.Range("F689").Value = "Number of participants with measurement for constipation following intervention"
$ws.Range("G689").Value = "GMHO:0000204"
$ws.Range("H689").Value = "measurement datum at followup"
$ws.Range("I689").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J689").Value = "measurement datum"
$ws.Range("K689").Value = "Intervention outcomes and spillover effects"
$ws.Range("L689").Value = "GMHO:0000178"
$ws.Range("M689").Value = "COMBO"
$ws.Range("N689").Value = "GMHO:0000257,GMHO:0000204"

# Row 690
$ws.Range("C690").Value = "1"
$ws.Range("E690").Value = "dizziness_n"
$ws.Range("F690").Value = "Number of participants with measurement for dizziness following intervention"
$ws.Range("G690").Value = "GMHO:0000179"
$ws.Range("H690").Value = "dizziness following intervention"
$ws.Range("I690").Value = "An adverse event following an intervention that involves experiencing dizziness."
$ws.Range("J690").Value = "adverse event following an intervention"
$ws.Range("K690").Value = "Intervention outcomes and spillover effects"
$ws.Range("L690").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M690").Value = "COMBO"
$ws.Range("N690").Value = "GMHO:0000179,GMHO:0000204"

# Row 691
$ws.Range("C691").Value = "1"
$ws.Range("E691").Value = "dizziness_n"
$ws.Range("F691").Value = "Number of participants with measurement for dizziness following intervention"
$ws.Range("G691").Value = "GMHO:0000206"
$ws.Range("H691").Value = "number of participants with measurement"
$ws.Range("I691").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J691").Value = "number of intervention participants"
$ws.Range("K691").Value = "Intervention outcomes and spillover effects"
$ws.Range("L691").Value = "GMHO:0000179,GMHO:0000204"
$ws.Range("M691").Value = "No Combo"

# Row 692
$ws.Range("C692").Value = "1"
$ws.Range("E692").Value = "dizziness_n"
$ws.Range("F692").Value = "Number of participants with measurement for dizziness following intervention"
$ws.Range("G692").Value = "GMHO:0000204"
$ws.Range("H692").Value = "measurement datum at followup"
$ws.Range("I692").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J692").Value = "measurement datum"
$ws.Range("K692").Value = "Intervention outcomes and spillover effects"
$ws.Range("L692").Value = "GMHO:0000047,GMHO:0000204"
$ws.Range("M692").Value = "COMBO"
$ws.Range("N692").Value = "GMHO:0000179,GMHO:0000204"

# Row 693
$ws.Range("C693").Value = "1"
$ws.Range("E693").Value = "vomiting_n"
$ws.Range("F693").Value = "Number of participants with measurement for vomiting following intervention"
$ws.Range("G693").Value = "GMHO:0000185"
$ws.Range("H693").Value = "vomiting following intervention"
$ws.Range("I693").Value = "A digestive system adverse event following intervention that has an outcome of vomiting, the retrograde expulsion of gastric contents through the oral cavity."
$ws.Range("J693").Value = "digestive system adverse event following intervention"
$ws.Range("K693").Value = "Intervention outcomes and spillover effects"
$ws.Range("L693").Value = "GMHO:0000178"
$ws.Range("M693").Value = "COMBO"
$ws.Range("N693").Value = "GMHO:0000185,GMHO:0000204"

# Row 694
$ws.Range("C694").Value = "1"
$ws.Range("E694").Value = "vomiting_n"
$ws.Range("F694").Value = "Number of participants with measurement for vomiting following intervention"
$ws.Range("G694").Value = "GMHO:0000206"
$ws.Range("H694").Value = "number of participants with measurement"
$ws.Range("I694").Value = "Number of intervention participants for whom a measurement was made."
$ws.Range("J694").Value = "number of intervention participants"
$ws.Range("K694").Value = "Intervention outcomes and spillover effects"
$ws.Range("L694").Value = "GMHO:0000185,GMHO:0000204"
$ws.Range("M694").Value = "No Combo"

# Row 695
$ws.Range("C695").Value = "1"
$ws.Range("E695").Value = "vomiting_n"
$ws.Range("F695").Value = "Number of participants with measurement for vomiting following intervention"
$ws.Range("G695").Value = "GMHO:0000204"
$ws.Range("H695").Value = "measurement datum at followup"
$ws.Range("I695").Value = "Measurement datum that was recorded as followup data in a study."
$ws.Range("J695").Value = "measurement datum"
$ws.Range("K695").Value = "Intervention outcomes and spillover effects"
$ws.Range("L695").Value = "GMHO:0000178"
$ws.Range("M695").Value = "COMBO"
$ws.Range("N695").Value = "GMHO:0000185,GMHO:0000204"
